$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D are stored as plain text in the workbook (they use
# "." as a thousands separator, e.g. "64.405.61"), so before writing any new
# value we force the cell to Text format. This stops Excel from silently
# re-interpreting numeric-looking strings (like "589.68") as real numbers,
# which would change the stored value (e.g. drop a trailing zero) and type.

# Apply the updated values cell by cell (row order follows the source data)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.405.61"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.502.73"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.68"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.51"
$ws.Range("E6").Value = "  +0.14%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.63"
$ws.Range("E9").Value = "  +5.92%  "
$ws.Range("E10").Value = "  +1.13%  "
$ws.Range("E11").Value = "  +4.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.103.54"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000181"
$ws.Range("E14").Value = "  +0.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.508.90"
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.78"
$ws.Range("E16").Value = "  +2.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.408.01"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("E18").Value = "  +0.60%  "
$ws.Range("E19").Value = "  +1.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.58"
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "390.80"
$ws.Range("E21").Value = "  +1.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.582"
$ws.Range("E22").Value = "  +3.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.646.05"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("E24").Value = "  +0.41%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000116"
$ws.Range("E27").Value = "  +3.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.45"
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("E30").Value = "  +1.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.21"
$ws.Range("E31").Value = "  -0.38%  "
$ws.Range("E32").Value = "  -4.52%  "
$ws.Range("E33").Value = "  +6.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.534.40"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.41"
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.34"
$ws.Range("E37").Value = "  +1.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.97"
$ws.Range("E38").Value = "  +1.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.56"
$ws.Range("E39").Value = "  +2.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "165.41"
$ws.Range("E40").Value = "  +2.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0791"
$ws.Range("E41").Value = "  +1.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.810"
$ws.Range("E42").Value = "  +0.63%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.45"
$ws.Range("E44").Value = "  +1.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "24.96"
$ws.Range("E45").Value = "  -2.40%  "
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("E47").Value = "  +1.04%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.425.75"
$ws.Range("E48").Value = "  -1.68%  "
$ws.Range("B49").Value = "SuiNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.928"
$ws.Range("E49").Value = "  +3.86%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.82"
$ws.Range("E50").Value = "  +1.36%  "
$ws.Range("E51").Value = "  +0.73%  "
